# Update dataSourceList sheet to match updated breakpoint definition schema:
# adds Contact/Notes columns, splits the old MTUSGS row into
# MTUSGSGlacier/MTUSGSYellowstone, adds a MEFRO row, and styles the table
# with a bordered header + body.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Write brand-new strings in the exact order needed so the shared
#        string table is rebuilt in the same sequence as the target file ---
$ws.Range("A9").Value  = "MTUSGSGlacier"
$ws.Range("A10").Value = "MTUSGSYellowstone"
$ws.Range("D1").Value  = "Contact"
$ws.Range("D10").Value = "Robert Al-Chokhachy"
$ws.Range("A8").Value  = "MEFRO"
$ws.Range("D7").Value  = "Jed Wright"
$ws.Range("E1").Value  = "Notes"
$ws.Range("E8").Value  = "Waiting for Daymet update to pair air temp"
$ws.Range("D4").Value  = "Ben Letcher"
$ws.Range("D9").Value  = "Leslie Jones"

# --- 2. Fill in the remaining (already-existing) string values ---
$ws.Range("B8").Value  = "ME"
$ws.Range("C8").Value  = "Northeast"
$ws.Range("D8").Value  = "Jed Wright"
$ws.Range("B9").Value  = "MT"
$ws.Range("C9").Value  = "West"
$ws.Range("B10").Value = "MT"
$ws.Range("C10").Value = "West"

# --- 3. Leave the remaining Contact/Notes body cells explicitly blank (but
#        present), so they still pick up the bordered style below ---
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("E10").Value = ""

# --- 4. Style: thin box border around every cell of the table, bold header ---
$ws.Range("A1:E10").Borders.LineStyle = 1
$ws.Range("A1:E1").Font.Bold = $true

# --- 5. Column widths for the new / widened columns ---
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null

# --- 6. Page setup + selection to match the saved view state ---
$ws.PageSetup.Orientation = 1
$ws.Range("D16").Select() | Out-Null

Write-Output "done"
